$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data rows 4 and 5: UserName -> phone number, Password -> "Tcs@1983" ---
$ws.Range("E4").Value = 9967887510
$ws.Range("F4").Value = "Tcs@1983"

$ws.Range("E5").Value = 9967887510
$ws.Range("F5").Value = "Tcs@1983"

# Add actual hyperlinks on F4 and F5 pointing at a mailto link for the password value.
$ws.Hyperlinks.Add($ws.Range("F4"), "mailto:Tcs@1983")
$ws.Hyperlinks.Add($ws.Range("F5"), "mailto:Tcs@1983")

# Copy the existing "Hyperlink" style (used by column D) onto F4/F5 so that the
# cell format matches a hyperlink-styled cell (style index 4 in styles.xml).
# (Hyperlinks.Add resets the cell style, so re-apply the format afterwards.)
$ws.Range("D2").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("F5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column E width change (10.42578125 -> 11) ---
$ws.Columns("E").ColumnWidth = 10.2

# --- Selection change ---
$ws.Range("E5:F5").Select()
